$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.949.11'
$ws.Range('E2').Value = '  +2.62%  '
$ws.Range('D3').Value = '3.199.16'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'536.20"
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = "'145.17"
$ws.Range('E6').Value = '  +3.93%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = "'0.531"
$ws.Range('E8').Value = '  +3.07%  '
$ws.Range('D9').Value = "'7.35"
$ws.Range('E9').Value = '  +0.34%  '
$ws.Range('E10').Value = '  +3.01%  '
$ws.Range('D11').Value = "'0.432"
$ws.Range('E11').Value = '  +2.57%  '
$ws.Range('D12').Value = '3.754.46'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('D14').Value = "'25.92"
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').Value = '60.028.11'
$ws.Range('E16').Value = '  +2.64%  '
$ws.Range('D17').Value = '3.197.26'
$ws.Range('E17').Value = '  +1.75%  '
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').Value = "'13.19"
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').Value = "'8.28"
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').Value = "'376.83"
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').Value = '  +1.79%  '
$ws.Range('D24').Value = "'70.07"
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').Value = "'0.169"
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = "'8.77"
$ws.Range('E26').Value = '  +7.75%  '
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('D28').Value = '0.0₃0894'
$ws.Range('E28').Value = '  +3.07%  '
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = "'22.34"
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').Value = "'6.16"
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').Value = "'5.42"
$ws.Range('E32').Value = '  +4.84%  '
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('D34').Value = "'6.66"
$ws.Range('E34').Value = '  +6.96%  '
$ws.Range('D35').Value = "'156.91"
$ws.Range('E35').Value = '  -2.35%  '
$ws.Range('D37').Value = '2.797.86'
$ws.Range('E37').Value = '  +5.69%  '
$ws.Range('D38').Value = "'25.54"
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('E39').Value = '  +3.37%  '
$ws.Range('E40').Value = '  +0.79%  '
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').Value = "'39.77"
$ws.Range('E42').Value = '  +2.80%  '
$ws.Range('E43').Value = '  +4.27%  '
$ws.Range('E44').Value = '  +1.52%  '
$ws.Range('E45').Value = '  +2.40%  '
$ws.Range('D46').Value = '3.243.14'
$ws.Range('E46').Value = '  +1.34%  '
$ws.Range('D47').Value = "'0.985"
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('E48').Value = '  +7.45%  '
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('D50').Value = "'20.58"
$ws.Range('E50').Value = '  +1.56%  '
$ws.Range('E51').Value = '  -0.02%  '
